# Update the "Förändrad" (Changed) date column (C) for all existing data
# rows (2-380) from 45190 to 45192.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C380").Value2 = 45192

# Row 380 gains an explicit custom row height (matches the rest of the sheet).
$ws.Rows.Item(380).RowHeight = 15

# Append the new record as row 381.
$ws.Range("A381").Value = "A 44627-2023"
$ws.Range("B381").Value2 = 45189
$ws.Range("C381").Value2 = 45192
$ws.Range("D381").Value = "SÖDERMANLANDS LÄN"
$ws.Range("E381").Value = "GNESTA"
$ws.Range("G381").Value2 = 1.9
$ws.Range("H381").Value2 = 0
$ws.Range("I381").Value2 = 0
$ws.Range("J381").Value2 = 0
$ws.Range("K381").Value2 = 0
$ws.Range("L381").Value2 = 0
$ws.Range("M381").Value2 = 0
$ws.Range("N381").Value2 = 0
$ws.Range("O381").Value2 = 0
$ws.Range("P381").Value2 = 0
$ws.Range("Q381").Value2 = 0
$ws.Range("R381").Value = ""

# Apply the same date format used elsewhere in columns B/C.
$ws.Range("B381:C381").NumberFormat = $ws.Range("B380:C380").NumberFormat

# Match the wrap-text style used in column R for the other rows.
$ws.Range("R381").WrapText = $ws.Range("R380").WrapText
